# Update specific numeric values on Sheet1 to match the new result data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C11"  = -12.8612
    "B12"  = 4.576899999999998
    "C23"  = -12.2648
    "D24"  = -8.036899999999992
    "C28"  = -13.7929
    "B32"  = 6.786499999999995
    "C32"  = -11.52660000000001
    "C34"  = -12.04500000000001
    "B36"  = 9.12310000000001
    "B38"  = 6.182700000000002
    "D38"  = -7.611400000000006
    "C42"  = -12.3714
    "B46"  = 6.276600000000003
    "D52"  = -7.687700000000004
    "B54"  = 4.999000000000004
    "C54"  = -13.94559999999999
    "B55"  = 5.626099999999997
    "B67"  = 5.517899999999995
    "B69"  = 5.353799999999996
    "B72"  = 5.138500000000006
    "D78"  = -7.565400000000002
    "D83"  = -9.079099999999999
    "D85"  = -8.846500000000001
    "D86"  = -8.619599999999997
    "B91"  = 4.748599999999993
    "D96"  = -8.284999999999997
    "C97"  = -11.5488
    "B99"  = 6.085800000000003
    "C99"  = -12.41690000000001
    "C101" = -12.91590000000001
    "D103" = -7.787900000000004
    "B104" = 9.841700000000005
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
